# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 397
$wsExhibit.Range("F5").Value = 19
$wsExhibit.Range("F7").Value = 0
$wsExhibit.Range("F8").Value = 0
$wsExhibit.Range("F9").Value = 65
$wsExhibit.Range("F10").Value = 524

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 397
$wsAll.Range("F4").Value = 1645
$wsAll.Range("F6").Value = 0
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 524
